# Auto-generated: updates market-price driven columns (H-N) on the
# Chocobo_Profits leve-profit sheets to match the scheduled runner's refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3996
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 4245
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 4245
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -4741
$ws.Range("H67").Value = 3996
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 4245
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 4245
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -5961
$ws.Range("H74").Value = 10006081
$ws.Range("I74").Value = 33337000
$ws.Range("J74").Value = 7116.4287
$ws.Range("K74").Value = 33337000
$ws.Range("L74").Value = 7116.4287
$ws.Range("M74").Value = -33336064
$ws.Range("N74").Value = -8988.4287
$ws.Range("H76").Value = 4642.7144
$ws.Range("I76").Value = 3999.8
$ws.Range("J76").Value = 6250
$ws.Range("K76").Value = 3999.8
$ws.Range("L76").Value = 6250
$ws.Range("M76").Value = -3684.8
$ws.Range("N76").Value = -6880
$ws.Range("H77").Value = 10006081
$ws.Range("I77").Value = 33337000
$ws.Range("J77").Value = 7116.4287
$ws.Range("K77").Value = 166685000
$ws.Range("L77").Value = 35582.14350000001
$ws.Range("M77").Value = -166680320
$ws.Range("N77").Value = -44942.14350000001
$ws.Range("H79").Value = 4642.7144
$ws.Range("I79").Value = 3999.8
$ws.Range("J79").Value = 6250
$ws.Range("K79").Value = 3999.8
$ws.Range("L79").Value = 6250
$ws.Range("M79").Value = -2907.8
$ws.Range("N79").Value = -8434
$ws.Range("H98").Value = 4468.3687
$ws.Range("I98").Value = 1699.8889
$ws.Range("J98").Value = 6960
$ws.Range("K98").Value = 1699.8889
$ws.Range("L98").Value = 6960
$ws.Range("M98").Value = -201.8888999999999
$ws.Range("N98").Value = -9956
$ws.Range("H115").Value = 1477.7778
$ws.Range("I115").Value = 1287.5
$ws.Range("K115").Value = 3862.5
$ws.Range("M115").Value = -2295.5
$ws.Range("H122").Value = 4468.3687
$ws.Range("I122").Value = 1699.8889
$ws.Range("J122").Value = 6960
$ws.Range("K122").Value = 5099.6667
$ws.Range("L122").Value = 20880
$ws.Range("M122").Value = -2649.6667
$ws.Range("N122").Value = -25780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 30000
$ws.Range("J18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30644
$ws.Range("H32").Value = 4725.5454
$ws.Range("I32").Value = 3841.7708
$ws.Range("K32").Value = 3841.7708
$ws.Range("M32").Value = -3554.7708
$ws.Range("H74").Value = 6392.8335
$ws.Range("I74").Value = 7598.1665
$ws.Range("K74").Value = 7598.1665
$ws.Range("M74").Value = -6724.1665
$ws.Range("H77").Value = 6392.8335
$ws.Range("I77").Value = 7598.1665
$ws.Range("K77").Value = 37990.8325
$ws.Range("M77").Value = -33622.8325
$ws.Range("H122").Value = 3096.4546
$ws.Range("I122").Value = 1260.75
$ws.Range("J122").Value = 7991.6665
$ws.Range("K122").Value = 3782.25
$ws.Range("L122").Value = 23974.9995
$ws.Range("M122").Value = -1332.25
$ws.Range("N122").Value = -28874.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2514500.2
$ws.Range("I7").Value = 4157.6
$ws.Range("J7").Value = 3655565
$ws.Range("K7").Value = 4157.6
$ws.Range("L7").Value = 3655565
$ws.Range("M7").Value = -4044.6
$ws.Range("N7").Value = -3655791
$ws.Range("H64").Value = 489.0625
$ws.Range("I64").Value = 587.3333
$ws.Range("J64").Value = 430.1
$ws.Range("K64").Value = 587.3333
$ws.Range("L64").Value = 430.1
$ws.Range("M64").Value = -362.3333
$ws.Range("N64").Value = -880.1
$ws.Range("H67").Value = 489.0625
$ws.Range("I67").Value = 587.3333
$ws.Range("J67").Value = 430.1
$ws.Range("K67").Value = 587.3333
$ws.Range("L67").Value = 430.1
$ws.Range("M67").Value = 192.6667
$ws.Range("N67").Value = -1990.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2794.75
$ws.Range("I5").Value = 166.66667
$ws.Range("J5").Value = 4371.6
$ws.Range("K5").Value = 166.66667
$ws.Range("L5").Value = 4371.6
$ws.Range("M5").Value = -54.66667000000001
$ws.Range("N5").Value = -4595.6
$ws.Range("H94").Value = 1283.5238
$ws.Range("I94").Value = 909.8
$ws.Range("J94").Value = 1623.2727
$ws.Range("K94").Value = 909.8
$ws.Range("L94").Value = 1623.2727
$ws.Range("M94").Value = -458.8
$ws.Range("N94").Value = -2525.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1507371.5
$ws.Range("J4").Value = 888.7778
$ws.Range("L4").Value = 2666.3334
$ws.Range("N4").Value = -2890.3334
$ws.Range("H80").Value = 17817.092
$ws.Range("I80").Value = 9900
$ws.Range("J80").Value = 18608.8
$ws.Range("K80").Value = 29700
$ws.Range("L80").Value = 55826.39999999999
$ws.Range("M80").Value = -28764
$ws.Range("N80").Value = -57698.39999999999
$ws.Range("H83").Value = 17817.092
$ws.Range("I83").Value = 9900
$ws.Range("J83").Value = 18608.8
$ws.Range("K83").Value = 89100
$ws.Range("L83").Value = 167479.2
$ws.Range("M83").Value = -84420
$ws.Range("N83").Value = -176839.2
$ws.Range("H92").Value = 713.9167
$ws.Range("J92").Value = 653.8570999999999
$ws.Range("L92").Value = 1961.5713
$ws.Range("N92").Value = -4457.5713
$ws.Range("H98").Value = 626
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H132").Value = 1975.8422
$ws.Range("I132").Value = 770.41174
$ws.Range("J132").Value = 2951.6667
$ws.Range("K132").Value = 6933.70566
$ws.Range("L132").Value = 26565.0003
$ws.Range("M132").Value = -4403.70566
$ws.Range("N132").Value = -31625.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8786.333000000001
$ws.Range("J5").Value = 12979.5
$ws.Range("L5").Value = 12979.5
$ws.Range("N5").Value = -13203.5
$ws.Range("H70").Value = 6625.2144
$ws.Range("I70").Value = 5873.1113
$ws.Range("J70").Value = 7979
$ws.Range("K70").Value = 5873.1113
$ws.Range("L70").Value = 7979
$ws.Range("M70").Value = -5603.1113
$ws.Range("N70").Value = -8519
$ws.Range("H73").Value = 6625.2144
$ws.Range("I73").Value = 5873.1113
$ws.Range("J73").Value = 7979
$ws.Range("K73").Value = 5873.1113
$ws.Range("L73").Value = 7979
$ws.Range("M73").Value = -4937.1113
$ws.Range("N73").Value = -9851
$ws.Range("H80").Value = 27780656
$ws.Range("I80").Value = 41669500
$ws.Range("J80").Value = 2966.6667
$ws.Range("K80").Value = 41669500
$ws.Range("L80").Value = 2966.6667
$ws.Range("M80").Value = -41668502
$ws.Range("N80").Value = -4962.6667
$ws.Range("H83").Value = 27780656
$ws.Range("I83").Value = 41669500
$ws.Range("J83").Value = 2966.6667
$ws.Range("K83").Value = 208347500
$ws.Range("L83").Value = 14833.3335
$ws.Range("M83").Value = -208342508
$ws.Range("N83").Value = -24817.3335
$ws.Range("H113").Value = 1532.2727
$ws.Range("I113").Value = 1581.25
$ws.Range("J113").Value = 1401.6666
$ws.Range("K113").Value = 1581.25
$ws.Range("L113").Value = 1401.6666
$ws.Range("M113").Value = 588.75
$ws.Range("N113").Value = -5741.6666
$ws.Range("H122").Value = 3864.7
$ws.Range("I122").Value = 1753
$ws.Range("J122").Value = 10199.8
$ws.Range("K122").Value = 5259
$ws.Range("L122").Value = 30599.4
$ws.Range("M122").Value = -2809
$ws.Range("N122").Value = -35499.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2573.1428
$ws.Range("J2").Value = 2573.1428
$ws.Range("L2").Value = 2573.1428
$ws.Range("N2").Value = -2797.1428
$ws.Range("H7").Value = 3333.0476
$ws.Range("I7").Value = 2612.6
$ws.Range("J7").Value = 5134.1665
$ws.Range("K7").Value = 2612.6
$ws.Range("L7").Value = 5134.1665
$ws.Range("M7").Value = -2500.6
$ws.Range("N7").Value = -5358.1665
$ws.Range("H50").Value = 49600
$ws.Range("J50").Value = 49600
$ws.Range("L50").Value = 49600
$ws.Range("N50").Value = -50874
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H126").Value = 3333.0476
$ws.Range("I126").Value = 2612.6
$ws.Range("J126").Value = 5134.1665
$ws.Range("K126").Value = 7837.799999999999
$ws.Range("L126").Value = 15402.4995
$ws.Range("M126").Value = -5367.799999999999
$ws.Range("N126").Value = -20342.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 7399.75
$ws.Range("J20").Value = 7399.75
$ws.Range("L20").Value = 7399.75
$ws.Range("N20").Value = -7879.75
$ws.Range("H113").Value = 547.6923
$ws.Range("I113").Value = 237.77777
$ws.Range("J113").Value = 1245
$ws.Range("K113").Value = 713.33331
$ws.Range("L113").Value = 3735
$ws.Range("M113").Value = 1456.66669
$ws.Range("N113").Value = -8075
$ws.Range("H122").Value = 2368.9666
$ws.Range("I122").Value = 1130.55
$ws.Range("J122").Value = 4845.8
$ws.Range("K122").Value = 3391.65
$ws.Range("L122").Value = 14537.4
$ws.Range("M122").Value = -941.6499999999996
$ws.Range("N122").Value = -19437.4
$ws.Range("H132").Value = 15875730
$ws.Range("I132").Value = 1882.1177
$ws.Range("K132").Value = 5646.3531
$ws.Range("M132").Value = -3116.3531
